$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "42÷7="
$t.Cell(1,2).Range.Text = "18÷8="
$t.Cell(1,3).Range.Text = "74÷2="
$t.Cell(1,4).Range.Text = "31÷3="
$t.Cell(1,5).Range.Text = "30÷4="
$t.Cell(5,1).Range.Text = "91÷9="
$t.Cell(5,2).Range.Text = "75÷3="
$t.Cell(5,3).Range.Text = "59÷6="
$t.Cell(5,4).Range.Text = "17÷8="
$t.Cell(5,5).Range.Text = "67÷3="
$t.Cell(9,1).Range.Text = "73÷6="
$t.Cell(9,2).Range.Text = "94÷6="
$t.Cell(9,3).Range.Text = "12÷8="
$t.Cell(9,4).Range.Text = "30÷5="
$t.Cell(9,5).Range.Text = "29÷8="
$t.Cell(13,1).Range.Text = "97÷2="
$t.Cell(13,2).Range.Text = "56÷9="
$t.Cell(13,3).Range.Text = "10÷4="
$t.Cell(13,4).Range.Text = "57÷4="
$t.Cell(13,5).Range.Text = "62÷9="
$t.Cell(17,1).Range.Text = "59÷5="
$t.Cell(17,2).Range.Text = "73÷4="
$t.Cell(17,3).Range.Text = "87÷7="
$t.Cell(17,4).Range.Text = "38÷3="
$t.Cell(17,5).Range.Text = "35÷6="
